$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Fix typo: "Praksa iz oblasti X" -> "Treninzi iz oblasti X"
# (commit: "Ispravka greske u pisanju" / "Fix a writing error")
# ---------------------------------------------------------------------------

# 1) "Praksa iz oblasti Blockchain" -> "Treninzi  iz oblasti Blockchain"
#    (kept as two runs with identical, non-bold formatting)
$rng = $d.Content
$found = $rng.Find.Execute("Praksa iz oblasti Blockchain", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $start = $rng.Start
    $rng.Text = "Treninzi  iz oblasti Blockchain"
    $split = $start + 9
    $r1 = $d.Range($start, $split)
    # toggling a property on the first piece and back forces Word to keep
    # it as a run boundary distinct from the following (identically
    # formatted) text, matching the two separate <w:r> runs in the target
    $r1.Bold = 1
    $r1.Bold = 0
}

# 2) "Praksa iz oblasti AI" -> "Treninzi" + " " (bold) + "iz oblasti AI"
$rng = $d.Content
$found = $rng.Find.Execute("Praksa iz oblasti AI", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $start = $rng.Start
    $rng.Text = "Treninzi iz oblasti AI"
    $midStart = $start + 8
    $midEnd = $start + 9
    $mid = $d.Range($midStart, $midEnd)
    $mid.Bold = 1
}

# 3) "Praksa iz oblasti DevOps" -> "Treninzi" + " " (bold) + "iz oblasti DevOps"
$rng = $d.Content
$found = $rng.Find.Execute("Praksa iz oblasti DevOps", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $start = $rng.Start
    $rng.Text = "Treninzi iz oblasti DevOps"
    $midStart = $start + 8
    $midEnd = $start + 9
    $mid = $d.Range($midStart, $midEnd)
    $mid.Bold = 1
}

Write-Output "done"
